# Update Ajinkya Rahane's Delhi Capitals per-innings batting activity
# (runs, balls, fours, sixes) to reflect data synced "till Excel form".
#
# The cells in this sheet store numeric-looking values as TEXT, so we
# force a text number format before writing each value and then reset
# the cell style back to Normal so no stray style/format is left behind
# (this keeps the cells as plain text values, matching the source data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Value
    )
    $rng = $ws.Range($Address)
    $rng.NumberFormat = "@"
    $rng.Value = $Value
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "C2" "60"
Set-TextValue "D2" "46"
Set-TextValue "E2" "5"
Set-TextValue "F2" "1"

# Row 3
Set-TextValue "C3" "26"
Set-TextValue "D3" "19"
Set-TextValue "E3" "3"
Set-TextValue "F3" "1"

# Row 4
Set-TextValue "C4" "0"
Set-TextValue "D4" "3"
Set-TextValue "E4" "0"
Set-TextValue "F4" "0"

# Row 5 unchanged

# Row 6
Set-TextValue "C6" "2"
Set-TextValue "D6" "9"
Set-TextValue "E6" "0"
Set-TextValue "F6" "0"

# Row 7
Set-TextValue "C7" "8"
Set-TextValue "D7" "10"
Set-TextValue "E7" "1"

# Row 8
Set-TextValue "C8" "15"
Set-TextValue "D8" "15"
Set-TextValue "E8" "3"

# Row 9
Set-TextValue "C9" "2"
Set-TextValue "D9" "4"
